$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.0593035
$ws.Range("H2").Value = 18.118607
$ws.Range("I2").Value = 0.2059870424264786
$ws.Range("J2").Value = 0.155107484181908
$ws.Range("M2").Value = 79.6373215
$ws.Range("N2").Value = 159.274643
$ws.Range("O2").Value = 0.6291200803678693
$ws.Range("P2").Value = 0.5421226517510304
$ws.Range("Q2").Value = 721.4586653955753
$ws.Range("R2").Value = 2885.834661582301
$ws.Range("S2").Value = 0.1295905846860859
$ws.Range("T2").Value = 0.08408728063112697

$ws.Range("G3").Value = 9.0593035
$ws.Range("H3").Value = 18.118607
$ws.Range("I3").Value = 0.2059870424264786
$ws.Range("J3").Value = 0.155107484181908
$ws.Range("O3").Value = 0.01453518499705013
$ws.Range("P3").Value = 0.01878779571703936
$ws.Range("Q3").Value = 16.66857488179
$ws.Range("R3").Value = 100.01144929074
$ws.Range("S3").Value = 0.00299405976866408
$ws.Range("T3").Value = 0.002914127726993601

$ws.Range("G4").Value = 9.0593035
$ws.Range("H4").Value = 18.118607
$ws.Range("I4").Value = 0.2059870424264786
$ws.Range("J4").Value = 0.155107484181908
$ws.Range("M4").Value = 19.69434466666667
$ws.Range("N4").Value = 59.083034
$ws.Range("O4").Value = 0.155581672840238
$ws.Range("P4").Value = 0.20110075566502
$ws.Range("Q4").Value = 178.4170455689397
$ws.Range("R4").Value = 1070.502273413638
$ws.Range("S4").Value = 0.03204780864412462
$ws.Range("T4").Value = 0.03119223227828183

$ws.Range("G5").Value = 9.0593035
$ws.Range("H5").Value = 18.118607
$ws.Range("I5").Value = 0.2059870424264786
$ws.Range("J5").Value = 0.155107484181908
$ws.Range("M5").Value = 6.320253
$ws.Range("N5").Value = 12.640506
$ws.Range("O5").Value = 0.04992882734391395
$ws.Range("P5").Value = 0.04302445450902571
$ws.Range("Q5").Value = 57.2570901237855
$ws.Range("R5").Value = 229.028360495142
$ws.Range("S5").Value = 0.01028469147639513
$ws.Range("T5").Value = 0.006673414897193926

$ws.Range("G6").Value = 9.0593035
$ws.Range("H6").Value = 18.118607
$ws.Range("I6").Value = 0.2059870424264786
$ws.Range("J6").Value = 0.155107484181908
$ws.Range("M6").Value = 6.942085
$ws.Range("N6").Value = 20.826255
$ws.Range("O6").Value = 0.05484118489746769
$ws.Range("P6").Value = 0.07088626522078065
$ws.Range("Q6").Value = 62.8904549377975
$ws.Range("R6").Value = 377.342729626785
$ws.Range("S6").Value = 0.01129657348019304
$ws.Range("T6").Value = 0.01099499026144677

$ws.Range("G7").Value = 9.0593035
$ws.Range("H7").Value = 18.118607
$ws.Range("I7").Value = 0.2059870424264786
$ws.Range("J7").Value = 0.155107484181908
$ws.Range("M7").Value = 12.151304
$ws.Range("N7").Value = 36.453912
$ws.Range("O7").Value = 0.09599304955346108
$ws.Range("P7").Value = 0.124078077137104
$ws.Range("Q7").Value = 110.082350856764
$ws.Range("R7").Value = 660.4941051405841
$ws.Range("S7").Value = 0.01977332437101586
$ws.Range("T7").Value = 0.01924543838686493

$ws.Range("G8").Value = 6.998093333333334
$ws.Range("I8").Value = 0.1591200193654833
$ws.Range("J8").Value = 0.1797251826815686
$ws.Range("M8").Value = 79.6373215
$ws.Range("N8").Value = 159.274643
$ws.Range("O8").Value = 0.6291200803678693
$ws.Range("P8").Value = 0.5421226517510304
$ws.Range("Q8").Value = 557.3094086736734
$ws.Range("R8").Value = 3343.85645204204
$ws.Range("S8").Value = 0.1001055993713498
$ws.Range("T8").Value = 0.09743309262177034

$ws.Range("G9").Value = 6.998093333333334
$ws.Range("I9").Value = 0.1591200193654833
$ws.Range("J9").Value = 0.1797251826815686
$ws.Range("O9").Value = 0.01453518499705013
$ws.Range("P9").Value = 0.01878779571703936
$ws.Range("Q9").Value = 12.87607184773334
$ws.Range("S9").Value = 0.002312838918211498
$ws.Range("T9").Value = 0.003376640017428891

$ws.Range("G10").Value = 6.998093333333334
$ws.Range("I10").Value = 0.1591200193654833
$ws.Range("J10").Value = 0.1797251826815686
$ws.Range("M10").Value = 19.69434466666667
$ws.Range("N10").Value = 59.083034
$ws.Range("O10").Value = 0.155581672840238
$ws.Range("P10").Value = 0.20110075566502
$ws.Range("Q10").Value = 137.8228621161689
$ws.Range("R10").Value = 1240.40575904552
$ws.Range("S10").Value = 0.02475615879525295
$ws.Range("T10").Value = 0.03614287004929721

$ws.Range("G11").Value = 6.998093333333334
$ws.Range("I11").Value = 0.1591200193654833
$ws.Range("J11").Value = 0.1797251826815686
$ws.Range("M11").Value = 6.320253
$ws.Range("N11").Value = 12.640506
$ws.Range("O11").Value = 0.04992882734391395
$ws.Range("P11").Value = 0.04302445450902571
$ws.Range("Q11").Value = 44.22972038428001
$ws.Range("R11").Value = 265.3783223056801
$ws.Range("S11").Value = 0.007944675973859462
$ws.Range("T11").Value = 0.007732577946409484

$ws.Range("G12").Value = 6.998093333333334
$ws.Range("I12").Value = 0.1591200193654833
$ws.Range("J12").Value = 0.1797251826815686
$ws.Range("M12").Value = 6.942085
$ws.Range("N12").Value = 20.826255
$ws.Range("O12").Value = 0.05484118489746769
$ws.Range("P12").Value = 0.07088626522078065
$ws.Range("Q12").Value = 48.58135875793334
$ws.Range("R12").Value = 437.2322288214
$ws.Range("S12").Value = 0.00872633040291111
$ws.Range("T12").Value = 0.01274004696641893

$ws.Range("G13").Value = 6.998093333333334
$ws.Range("I13").Value = 0.1591200193654833
$ws.Range("J13").Value = 0.1797251826815686
$ws.Range("M13").Value = 12.151304
$ws.Range("N13").Value = 36.453912
$ws.Range("O13").Value = 0.09599304955346108
$ws.Range("P13").Value = 0.124078077137104
$ws.Range("Q13").Value = 85.03595951370669
$ws.Range("R13").Value = 765.3236356233601
$ws.Range("S13").Value = 0.01527441590389853
$ws.Range("T13").Value = 0.02229995508024379

$ws.Range("G14").Value = 8.095122333333334
$ws.Range("H14").Value = 24.285367
$ws.Range("I14").Value = 0.1840638529798531
$ws.Range("J14").Value = 0.2078991049259102
$ws.Range("M14").Value = 79.6373215
$ws.Range("N14").Value = 159.274643
$ws.Range("O14").Value = 0.6291200803678693
$ws.Range("P14").Value = 0.5421226517510304
$ws.Range("Q14").Value = 644.6738598414969
$ws.Range("R14").Value = 3868.043159048981
$ws.Range("S14").Value = 0.1157982659795048
$ws.Range("T14").Value = 0.1127068140591001

$ws.Range("G15").Value = 8.095122333333334
$ws.Range("H15").Value = 24.285367
$ws.Range("I15").Value = 0.1840638529798531
$ws.Range("J15").Value = 0.2078991049259102
$ws.Range("O15").Value = 0.01453518499705013
$ws.Range("P15").Value = 0.01878779571703936
$ws.Range("Q15").Value = 14.89453938599334
$ws.Range("R15").Value = 134.05085447394
$ws.Range("S15").Value = 0.002675402154332
$ws.Range("T15").Value = 0.003905965913103332

$ws.Range("G16").Value = 8.095122333333334
$ws.Range("H16").Value = 24.285367
$ws.Range("I16").Value = 0.1840638529798531
$ws.Range("J16").Value = 0.2078991049259102
$ws.Range("M16").Value = 19.69434466666667
$ws.Range("N16").Value = 59.083034
$ws.Range("O16").Value = 0.155581672840238
$ws.Range("P16").Value = 0.20110075566502
$ws.Range("Q16").Value = 159.4281293514976
$ws.Range("R16").Value = 1434.853164163478
$ws.Range("S16").Value = 0.02863696215602516
$ws.Range("T16").Value = 0.04180866710268182

$ws.Range("G17").Value = 8.095122333333334
$ws.Range("H17").Value = 24.285367
$ws.Range("I17").Value = 0.1840638529798531
$ws.Range("J17").Value = 0.2078991049259102
$ws.Range("M17").Value = 6.320253
$ws.Range("N17").Value = 12.640506
$ws.Range("O17").Value = 0.04992882734391395
$ws.Range("P17").Value = 0.04302445450902571
$ws.Range("Q17").Value = 51.163221212617
$ws.Range("R17").Value = 306.979327275702
$ws.Range("S17").Value = 0.009190092335686646
$ws.Range("T17").Value = 0.008944745582351985

$ws.Range("G18").Value = 8.095122333333334
$ws.Range("H18").Value = 24.285367
$ws.Range("I18").Value = 0.1840638529798531
$ws.Range("J18").Value = 0.2078991049259102
$ws.Range("M18").Value = 6.942085
$ws.Range("N18").Value = 20.826255
$ws.Range("O18").Value = 0.05484118489746769
$ws.Range("P18").Value = 0.07088626522078065
$ws.Range("Q18").Value = 56.19702732339834
$ws.Range("R18").Value = 505.773245910585
$ws.Range("S18").Value = 0.01009427979420843
$ws.Range("T18").Value = 0.01473719109094097

$ws.Range("G19").Value = 8.095122333333334
$ws.Range("H19").Value = 24.285367
$ws.Range("I19").Value = 0.1840638529798531
$ws.Range("J19").Value = 0.2078991049259102
$ws.Range("M19").Value = 12.151304
$ws.Range("N19").Value = 36.453912
$ws.Range("O19").Value = 0.09599304955346108
$ws.Range("P19").Value = 0.124078077137104
$ws.Range("Q19").Value = 98.36629238952268
$ws.Range("R19").Value = 885.2966315057041
$ws.Range("S19").Value = 0.01766885056009601
$ws.Range("T19").Value = 0.02579572117773197

$ws.Range("G20").Value = 6.067365000000001
$ws.Range("H20").Value = 12.13473
$ws.Range("I20").Value = 0.1379574678861274
$ws.Range("J20").Value = 0.1038814651439112
$ws.Range("M20").Value = 79.6373215
$ws.Range("N20").Value = 159.274643
$ws.Range("O20").Value = 0.6291200803678693
$ws.Range("P20").Value = 0.5421226517510304
$ws.Range("Q20").Value = 483.1886971628475
$ws.Range("R20").Value = 1932.75478865139
$ws.Range("S20").Value = 0.0867918132838682
$ws.Range("T20").Value = 0.05631649535159935

$ws.Range("G21").Value = 6.067365000000001
$ws.Range("H21").Value = 12.13473
$ws.Range("I21").Value = 0.1379574678861274
$ws.Range("J21").Value = 0.1038814651439112
$ws.Range("O21").Value = 0.01453518499705013
$ws.Range("P21").Value = 0.01878779571703936
$ws.Range("Q21").Value = 11.1635875581
$ws.Range("R21").Value = 66.9815253486
$ws.Range("S21").Value = 0.002005237317449463
$ws.Range("T21").Value = 0.001951703745910547

$ws.Range("G22").Value = 6.067365000000001
$ws.Range("H22").Value = 12.13473
$ws.Range("I22").Value = 0.1379574678861274
$ws.Range("J22").Value = 0.1038814651439112
$ws.Range("M22").Value = 19.69434466666667
$ws.Range("N22").Value = 59.083034
$ws.Range("O22").Value = 0.155581672840238
$ws.Range("P22").Value = 0.20110075566502
$ws.Range("Q22").Value = 119.49277752847
$ws.Range("R22").Value = 716.9566651708201
$ws.Range("S22").Value = 0.02146365363452711
$ws.Range("T22").Value = 0.02089064114002996

$ws.Range("G23").Value = 6.067365000000001
$ws.Range("H23").Value = 12.13473
$ws.Range("I23").Value = 0.1379574678861274
$ws.Range("J23").Value = 0.1038814651439112
$ws.Range("M23").Value = 6.320253
$ws.Range("N23").Value = 12.640506
$ws.Range("O23").Value = 0.04992882734391395
$ws.Range("P23").Value = 0.04302445450902571
$ws.Range("Q23").Value = 38.347281843345
$ws.Range("R23").Value = 153.38912737338
$ws.Range("S23").Value = 0.006888054594890009
$ws.Range("T23").Value = 0.004469443371415145

$ws.Range("G24").Value = 6.067365000000001
$ws.Range("H24").Value = 12.13473
$ws.Range("I24").Value = 0.1379574678861274
$ws.Range("J24").Value = 0.1038814651439112
$ws.Range("M24").Value = 6.942085
$ws.Range("N24").Value = 20.826255
$ws.Range("O24").Value = 0.05484118489746769
$ws.Range("P24").Value = 0.07088626522078065
$ws.Range("Q24").Value = 42.120163556025
$ws.Range("R24").Value = 252.72098133615
$ws.Range("S24").Value = 0.007565751004329574
$ws.Range("T24").Value = 0.007363769089714567

$ws.Range("G25").Value = 6.067365000000001
$ws.Range("H25").Value = 12.13473
$ws.Range("I25").Value = 0.1379574678861274
$ws.Range("J25").Value = 0.1038814651439112
$ws.Range("M25").Value = 12.151304
$ws.Range("N25").Value = 36.453912
$ws.Range("O25").Value = 0.09599304955346108
$ws.Range("P25").Value = 0.124078077137104
$ws.Range("Q25").Value = 73.72639659396002
$ws.Range("R25").Value = 442.3583795637601
$ws.Range("S25").Value = 0.01324295805106304
$ws.Range("T25").Value = 0.01288941244524159

$ws.Range("G26").Value = 5.882126666666667
$ws.Range("H26").Value = 17.64638
$ws.Range("I26").Value = 0.1337455881950073
$ws.Range("J26").Value = 0.1510649028768016
$ws.Range("M26").Value = 79.6373215
$ws.Range("N26").Value = 159.274643
$ws.Range("O26").Value = 0.6291200803678693
$ws.Range("P26").Value = 0.5421226517510304
$ws.Range("Q26").Value = 468.4368124570567
$ws.Range("R26").Value = 2810.62087474234
$ws.Range("S26").Value = 0.08414203519409091
$ws.Range("T26").Value = 0.08189570573408354

$ws.Range("G27").Value = 5.882126666666667
$ws.Range("H27").Value = 17.64638
$ws.Range("I27").Value = 0.1337455881950073
$ws.Range("J27").Value = 0.1510649028768016
$ws.Range("O27").Value = 0.01453518499705013
$ws.Range("P27").Value = 0.01878779571703936
$ws.Range("Q27").Value = 10.82276013906667
$ws.Range("R27").Value = 97.4048412516
$ws.Range("S27").Value = 0.001944016866953714
$ws.Range("T27").Value = 0.00283817653526374

$ws.Range("G28").Value = 5.882126666666667
$ws.Range("H28").Value = 17.64638
$ws.Range("I28").Value = 0.1337455881950073
$ws.Range("J28").Value = 0.1510649028768016
$ws.Range("M28").Value = 19.69434466666667
$ws.Range("N28").Value = 59.083034
$ws.Range("O28").Value = 0.155581672840238
$ws.Range("P28").Value = 0.20110075566502
$ws.Range("Q28").Value = 115.8446299463244
$ws.Range("R28").Value = 1042.60166951692
$ws.Range("S28").Value = 0.02080836234638081
$ws.Range("T28").Value = 0.03037926612298765

$ws.Range("G29").Value = 5.882126666666667
$ws.Range("H29").Value = 17.64638
$ws.Range("I29").Value = 0.1337455881950073
$ws.Range("J29").Value = 0.1510649028768016
$ws.Range("M29").Value = 6.320253
$ws.Range("N29").Value = 12.640506
$ws.Range("O29").Value = 0.04992882734391395
$ws.Range("P29").Value = 0.04302445450902571
$ws.Range("Q29").Value = 37.17652871138
$ws.Range("R29").Value = 223.05917226828
$ws.Range("S29").Value = 0.006677760380998734
$ws.Range("T29").Value = 0.006499485041733337

$ws.Range("G30").Value = 5.882126666666667
$ws.Range("H30").Value = 17.64638
$ws.Range("I30").Value = 0.1337455881950073
$ws.Range("J30").Value = 0.1510649028768016
$ws.Range("M30").Value = 6.942085
$ws.Range("N30").Value = 20.826255
$ws.Range("O30").Value = 0.05484118489746769
$ws.Range("P30").Value = 0.07088626522078065
$ws.Range("Q30").Value = 40.83422330076667
$ws.Range("R30").Value = 367.5080097069
$ws.Range("S30").Value = 0.007334766531422965
$ws.Range("T30").Value = 0.01070842677087643

$ws.Range("G31").Value = 5.882126666666667
$ws.Range("H31").Value = 17.64638
$ws.Range("I31").Value = 0.1337455881950073
$ws.Range("J31").Value = 0.1510649028768016
$ws.Range("M31").Value = 12.151304
$ws.Range("N31").Value = 36.453912
$ws.Range("O31").Value = 0.09599304955346108
$ws.Range("P31").Value = 0.124078077137104
$ws.Range("Q31").Value = 71.47550929317335
$ws.Range("R31").Value = 643.2795836385601
$ws.Range("S31").Value = 0.01283864687516013
$ws.Range("T31").Value = 0.01874384267185692

$ws.Range("G32").Value = 7.877956999999999
$ws.Range("H32").Value = 23.633871
$ws.Range("I32").Value = 0.1791260291470502
$ws.Range("J32").Value = 0.2023218601899006
$ws.Range("M32").Value = 79.6373215
$ws.Range("N32").Value = 159.274643
$ws.Range("O32").Value = 0.6291200803678693
$ws.Range("P32").Value = 0.5421226517510304
$ws.Range("Q32").Value = 627.3793943721754
$ws.Range("R32").Value = 3764.276366233053
$ws.Range("S32").Value = 0.1126917818529695
$ws.Range("T32").Value = 0.1096832633533501

$ws.Range("G33").Value = 7.877956999999999
$ws.Range("H33").Value = 23.633871
$ws.Range("I33").Value = 0.1791260291470502
$ws.Range("J33").Value = 0.2023218601899006
$ws.Range("O33").Value = 0.01453518499705013
$ws.Range("P33").Value = 0.01878779571703936
$ws.Range("Q33").Value = 14.49496820258
$ws.Range("R33").Value = 130.45471382322
$ws.Range("S33").Value = 0.002603629971439368
$ws.Range("T33").Value = 0.00380118177833925

$ws.Range("G34").Value = 7.877956999999999
$ws.Range("H34").Value = 23.633871
$ws.Range("I34").Value = 0.1791260291470502
$ws.Range("J34").Value = 0.2023218601899006
$ws.Range("M34").Value = 19.69434466666667
$ws.Range("N34").Value = 59.083034
$ws.Range("O34").Value = 0.155581672840238
$ws.Range("P34").Value = 0.20110075566502
$ws.Range("Q34").Value = 155.1512004271793
$ws.Range("R34").Value = 1396.360803844614
$ws.Range("S34").Value = 0.0278687272639273
$ws.Range("T34").Value = 0.04068707897174153

$ws.Range("G35").Value = 7.877956999999999
$ws.Range("H35").Value = 23.633871
$ws.Range("I35").Value = 0.1791260291470502
$ws.Range("J35").Value = 0.2023218601899006
$ws.Range("M35").Value = 6.320253
$ws.Range("N35").Value = 12.640506
$ws.Range("O35").Value = 0.04992882734391395
$ws.Range("P35").Value = 0.04302445450902571
$ws.Range("Q35").Value = 49.790681363121
$ws.Range("R35").Value = 298.744088178726
$ws.Range("S35").Value = 0.00894355258208397
$ws.Range("T35").Value = 0.008704787669921837

$ws.Range("G36").Value = 7.877956999999999
$ws.Range("H36").Value = 23.633871
$ws.Range("I36").Value = 0.1791260291470502
$ws.Range("J36").Value = 0.2023218601899006
$ws.Range("M36").Value = 6.942085
$ws.Range("N36").Value = 20.826255
$ws.Range("O36").Value = 0.05484118489746769
$ws.Range("P36").Value = 0.07088626522078065
$ws.Range("Q36").Value = 54.68944712034499
$ws.Range("R36").Value = 492.205024083105
$ws.Range("S36").Value = 0.009823483684402569
$ws.Range("T36").Value = 0.014341841041383

$ws.Range("G37").Value = 7.877956999999999
$ws.Range("H37").Value = 23.633871
$ws.Range("I37").Value = 0.1791260291470502
$ws.Range("J37").Value = 0.2023218601899006
$ws.Range("M37").Value = 12.151304
$ws.Range("N37").Value = 36.453912
$ws.Range("O37").Value = 0.09599304955346108
$ws.Range("P37").Value = 0.124078077137104
$ws.Range("Q37").Value = 95.72745040592801
$ws.Range("R37").Value = 861.547053653352
$ws.Range("S37").Value = 0.01719485379222751
$ws.Range("T37").Value = 0.02510370737516486
